$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename sheet tab
$ws.Name = "Seeds for Brands Collection"

# 2. Update row 11 (previously the "typo in both fields" test case data
#    was split across rows 11 and 12; the two rows are merged into one
#    fully-populated row 11). Shared strings get appended in this order
#    so that new-string indices line up with the target workbook.
$ws.Range("C11").Value = "typo in field name i.e. brandname`n / brandsName / .."
$ws.Range("C11").HorizontalAlignment = -4108
$ws.Range("C11").WrapText = $true

$ws.Range("D11").Value = "InvalidYear ; i.e. 2030/ 999 `n/ undefined / null / .."
$ws.Range("D11").HorizontalAlignment = -4108
$ws.Range("D11").WrapText = $true

$ws.Range("E11").Value = "InvalidNumber ; i.e. -1/ 0 `n/ undefined / null / .."
$ws.Range("E11").HorizontalAlignment = -4108
$ws.Range("E11").WrapText = $true

$ws.Range("F11").Value = "typo in field name i.e. headquarter`n / headQuarters / headAddress / hqAddress / .."
$ws.Range("F11").HorizontalAlignment = -4108
$ws.Range("F11").WrapText = $true

$ws.Range("G11").Value = "Fix all typeos in both fields `nand set the min value for both invalid field `nand Insert Successfully "
$ws.Range("G11").HorizontalAlignment = -4108
$ws.Range("G11").WrapText = $true

$ws.Range("H11").Value = "brandName,`nyearFounded , `nnumberOfLocations,`nheadquarters"
$ws.Range("H11").HorizontalAlignment = -4108
$ws.Range("H11").WrapText = $true

# 3. Rename the last header (H1) -- this also renames the Table2 ListColumn
#    "Field for TestCase" -> " Violated Field(s) by TestCase"
$ws.Range("H1").Value = " Violated Field(s) by TestCase"

# 4. Remove the now-redundant old row 12 (its data was folded into row 11
#    above), shrinking the table/sheet from 12 to 11 rows.
$ws.Rows.Item(12).Delete()

# 5. The merged row now holds 4 lines of wrapped text, so it is taller.
$ws.Rows.Item(11).RowHeight = 60

# 6. Column H needs to be a bit wider for the new header text.
$ws.Columns.Item(8).ColumnWidth = 29.140625

# 7. Keep the page in portrait orientation (explicit page setup).
$ws.PageSetup.Orientation = 1

# 8. Restore the active cell selection to G11 (the last edited cell).
[void]$ws.Range("G11").Select()

Write-Output "Workbook updated"
